$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores Price/Volume columns as literal text (not numbers),
# e.g. "64.462.03" or "  -0.76%  ". Before writing new values into the cells
# that change, force those specific cells to Text format ("@") so Excel does
# not reinterpret numeric-looking strings (e.g. "6.12", "0.568") as numbers.

$ws.Range("D2:E2").NumberFormat = "@"
# Row 2
$ws.Range("D2").Value = "63.669.55"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3:E3").NumberFormat = "@"
# Row 3
$ws.Range("D3").Value = "3.411.54"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("E4").NumberFormat = "@"
# Row 4
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5:E5").NumberFormat = "@"
# Row 5
$ws.Range("D5").Value = "566.84"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("D6:E6").NumberFormat = "@"
# Row 6
$ws.Range("D6").Value = "156.52"
$ws.Range("E6").Value = "  -2.15%  "

$ws.Range("E7").NumberFormat = "@"
# Row 7
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8:E8").NumberFormat = "@"
# Row 8
$ws.Range("D8").Value = "3.412.62"
$ws.Range("E8").Value = "  -1.40%  "

$ws.Range("D9:E9").NumberFormat = "@"
# Row 9
$ws.Range("D9").Value = "0.568"
$ws.Range("E9").Value = "  -6.74%  "

$ws.Range("D10:E10").NumberFormat = "@"
# Row 10
$ws.Range("D10").Value = "7.18"
$ws.Range("E10").Value = "  -1.38%  "

$ws.Range("D11:E11").NumberFormat = "@"
# Row 11
$ws.Range("D11").Value = "0.119"
$ws.Range("E11").Value = "  -4.86%  "

$ws.Range("D12:E12").NumberFormat = "@"
# Row 12
$ws.Range("D12").Value = "0.427"
$ws.Range("E12").Value = "  -5.05%  "

$ws.Range("D13:E13").NumberFormat = "@"
# Row 13
$ws.Range("D13").Value = "4.016.26"
$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("E14").NumberFormat = "@"
# Row 14
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D15:E15").NumberFormat = "@"
# Row 15
$ws.Range("D15").Value = "27.04"
$ws.Range("E15").Value = "  -4.44%  "

$ws.Range("D16:E16").NumberFormat = "@"
# Row 16
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -10.02%  "

$ws.Range("D17:E17").NumberFormat = "@"
# Row 17
$ws.Range("D17").Value = "63.782.74"
$ws.Range("E17").Value = "  -1.95%  "

$ws.Range("D18:E18").NumberFormat = "@"
# Row 18
$ws.Range("D18").Value = "3.457.03"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("D19:E19").NumberFormat = "@"
# Row 19
$ws.Range("D19").Value = "6.12"
$ws.Range("E19").Value = "  -5.29%  "

$ws.Range("D20:E20").NumberFormat = "@"
# Row 20
$ws.Range("D20").Value = "13.61"
$ws.Range("E20").Value = "  -4.70%  "

$ws.Range("D21:E21").NumberFormat = "@"
# Row 21
$ws.Range("D21").Value = "374.30"
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22:E22").NumberFormat = "@"
# Row 22
$ws.Range("D22").Value = "7.75"
$ws.Range("E22").Value = "  -4.87%  "

$ws.Range("D23:E23").NumberFormat = "@"
# Row 23
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("D24:E24").NumberFormat = "@"
# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "71.82"
$ws.Range("E24").Value = "  -1.56%  "

$ws.Range("D25:E25").NumberFormat = "@"
# Row 25
$ws.Range("D25").Value = "0.521"
$ws.Range("E25").Value = "  -6.54%  "

$ws.Range("D26:E26").NumberFormat = "@"
# Row 26
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  -2.30%  "

$ws.Range("D27:E27").NumberFormat = "@"
# Row 27
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -5.21%  "

$ws.Range("D28:E28").NumberFormat = "@"
# Row 28
$ws.Range("D28").Value = "0.177"
$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("E29").NumberFormat = "@"
# Row 29
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30:E30").NumberFormat = "@"
# Row 30
$ws.Range("D30").Value = "5.97"
$ws.Range("E30").Value = "  -2.92%  "

$ws.Range("D31:E31").NumberFormat = "@"
# Row 31
$ws.Range("D31").Value = "1.39"
$ws.Range("E31").Value = "  -8.85%  "

$ws.Range("D32:E32").NumberFormat = "@"
# Row 32
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -2.92%  "

$ws.Range("D33:E33").NumberFormat = "@"
# Row 33
$ws.Range("D33").Value = "22.93"
$ws.Range("E33").Value = "  -2.84%  "

$ws.Range("D34:E34").NumberFormat = "@"
# Row 34
$ws.Range("D34").Value = "6.92"
$ws.Range("E34").Value = "  -5.51%  "

$ws.Range("D35:E35").NumberFormat = "@"
# Row 35
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  -6.03%  "

$ws.Range("D36:E36").NumberFormat = "@"
# Row 36
$ws.Range("D36").Value = "160.22"
$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("D37:E37").NumberFormat = "@"
# Row 37
$ws.Range("D37").Value = "1.83"
$ws.Range("E37").Value = "  -5.45%  "

$ws.Range("D38:E38").NumberFormat = "@"
# Row 38
$ws.Range("D38").Value = "0.820"
$ws.Range("E38").Value = "  +4.95%  "

$ws.Range("D39:E39").NumberFormat = "@"
# Row 39
$ws.Range("D39").Value = "26.36"
$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("D40:E40").NumberFormat = "@"
# Row 40
$ws.Range("D40").Value = "0.0731"
$ws.Range("E40").Value = "  -7.05%  "

$ws.Range("D41:E41").NumberFormat = "@"
# Row 41
$ws.Range("D41").Value = "2.779.22"
$ws.Range("E41").Value = "  -4.71%  "

$ws.Range("D42:E42").NumberFormat = "@"
# Row 42
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "42.56"
$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("D43:E43").NumberFormat = "@"
# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "4.42"
$ws.Range("E43").Value = "  -7.19%  "

$ws.Range("D44:E44").NumberFormat = "@"
# Row 44
$ws.Range("D44").Value = "6.33"
$ws.Range("E44").Value = "  -7.69%  "

$ws.Range("D45:E45").NumberFormat = "@"
# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0304"
$ws.Range("E45").Value = "  -4.94%  "

$ws.Range("D46:E46").NumberFormat = "@"
# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "25.15"
$ws.Range("E46").Value = "  -3.08%  "

$ws.Range("D47:E47").NumberFormat = "@"
# Row 47
$ws.Range("D47").Value = "2.33"
$ws.Range("E47").Value = "  +6.45%  "

$ws.Range("D48:E48").NumberFormat = "@"
# Row 48
$ws.Range("D48").Value = "326.17"
$ws.Range("E48").Value = "  +1.32%  "

$ws.Range("D49:E49").NumberFormat = "@"
# Row 49
$ws.Range("D49").Value = "1.03"
$ws.Range("E49").Value = "  -5.57%  "

$ws.Range("D50:E50").NumberFormat = "@"
# Row 50
$ws.Range("D50").Value = "6.34"
$ws.Range("E50").Value = "  -4.18%  "

$ws.Range("D51:E51").NumberFormat = "@"
# Row 51
$ws.Range("D51").Value = "0.823"
$ws.Range("E51").Value = "  -6.40%  "
